$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New values introduced by this commit -------------------------------
# A new "X-1P-ZUUL-HOST" header/value used in column F for every test row,
# and a new "PASS" expected-result value used in the new column L.
$hostHeaderValue = "X-1P-ZUUL-HOST=http://ec2-54-148-94-143.us-west-2.compute.amazonaws.com:7001/"
$passValue = "PASS"

# Row 2 (S1_TC_T1)
$ws.Range("F2").Value = $hostHeaderValue
$ws.Range("L2").Value = $passValue

# Row 3 (S1_TC_T2)
$ws.Range("F3").Value = $hostHeaderValue
$ws.Range("L3").Value = $passValue

# Row 4 (S1_TC_T3)
$ws.Range("F4").Value = $hostHeaderValue
$ws.Range("L4").Value = $passValue

# Row 2 height shrinks now that the long validation text no longer needs
# as much vertical room (45 -> 30, matching rows 3 & 4).
$ws.Rows.Item(2).RowHeight = 30

# Column F needs to widen considerably to fit the new long header text.
$ws.Columns.Item(6).ColumnWidth = 78.16666666666667

# The active selection when the workbook was saved moved to F10.
$ws.Range("F10").Select()
